$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75, shifting existing rows 75-186 down to 76-187.
$ws.Rows("75").Insert()

# Populate the newly inserted row 75 with the new data record.
$ws.Range("A75").Value = 11
$ws.Range("B75").Value = "Vega Monumental Concepción"
$ws.Range("C75").Value = "Bíobío"
$ws.Range("D75").Value = 45175
$ws.Range("E75").Value = 8
$ws.Range("F75").Value = 100112001
$ws.Range("G75").Value = "Berenjena"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 50
$ws.Range("K75").Value = 10000
$ws.Range("L75").Value = 10000
$ws.Range("M75").Value = 10000
$ws.Range("N75").Value = "$/caja 50 unidades"
$ws.Range("O75").Value = "Región de Arica y Parinacota"
$ws.Range("P75").Value = 200
$ws.Range("Q75").Value = 50
$ws.Range("R75").Value = "Hortaliza"
